# Daily 100 Error Counts.xlsx — "Add files via upload"
#
# New day of data (2025-10-23 / serial 45953) was appended to the table in
# row 13, and the in-sheet selection left sitting on cell J12 (outside the
# data block) when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 was a blank placeholder row (just the date-formatted A13 cell);
# fill it in with the next day's counts, continuing the existing table.
$ws.Range("A13").Value = 45953
$ws.Range("B13").Value = 686
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 656

# Selection moved off the table to J12 by the time the workbook was saved.
[void]$ws.Range("J12").Select()
